# Auto-generated: apply updated cryptocurrency price/volume snapshot
# (GitHub Actions refresh of cryptos list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '96.463.49'
$ws.Range('E2').Value = '  -0.37%  '
# Row 3
$ws.Range('D3').Value = '3.703.80'
$ws.Range('E3').Value = '  -0.17%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.33'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.05%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.87'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.53%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '650.34'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.55%  '
# Row 8
$ws.Range('E8').Value = '  -0.32%  '
# Row 9
$ws.Range('E9').Value = '  +0.00%  '
# Row 10
$ws.Range('E10').Value = '  -5.39%  '
# Row 11
$ws.Range('D11').Value = '3.704.34'
$ws.Range('E11').Value = '  -0.09%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000310'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +18.13%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.25'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.69%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.205'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.06%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.73'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.90%  '
# Row 16
$ws.Range('D16').Value = '4.390.59'
$ws.Range('E16').Value = '  -0.25%  '
# Row 17
$ws.Range('D17').Value = '96.240.18'
$ws.Range('E17').Value = '  -0.37%  '
# Row 18
$ws.Range('E18').Value = '  +12.84%  '
# Row 19
$ws.Range('D19').Value = '3.698.77'
$ws.Range('E19').Value = '  -0.13%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.05'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.62%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.67'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.58%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.501'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -7.01%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '519.92'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.17%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.40'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.88%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000210'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.45%  '
# Row 26
$ws.Range('E26').Value = '  +0.33%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '100.95'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.57%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.18'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.33%  '
# Row 29
$ws.Range('E29').Value = '  +2.48%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.01'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.69%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.12'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.57%  '
# Row 32
$ws.Range('E32').Value = '  +0.09%  '
# Row 33
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.187'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.73%  '
# Row 34
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +7.29%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.10%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.17'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.83%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '648.08'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.32%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.586'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.44%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.81'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.88%  '
# Row 40
$ws.Range('E40').Value = '  +0.03%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +11.82%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.04'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.66%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.55'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.83%  '
# Row 44
$ws.Range('E44').Value = '  +0.16%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.957'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.54%  '
# Row 46
$ws.Range('E46').Value = '  +1.80%  '
# Row 47
$ws.Range('E47').Value = '  +2.05%  '
# Row 48
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.28'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.21%  '
# Row 49
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.57'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.07%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.46'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.50%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.52'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.03%  '
